$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.992.14'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.823.25'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.49%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4686'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3670'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.68%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07362'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8749'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.30'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.843.33'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.431'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07225'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.524'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.80'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.63%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.004'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008757'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.57%  '
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.69'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.002.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.293'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.62'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.049.63'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.54%  '
$ws.Range("E25").Value = '  -0.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.94'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.147'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.245'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08861'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7556'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.163'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.510'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.934'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.24%  '
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.097'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05313'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01951'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.979'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.64%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.382'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.197'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5307'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1654'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.486'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4900'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.45'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.54%  '
$ws.Range("E48").Value = '  +0.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.665'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.12'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06298'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.30%  '
